$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Jídlo 1 Nevim co už dál / Traditional Hungarian goulash with pork meat and sauerkraut served with bread dumplings"
$ws.Range("C3").Value = "Jídlo 2 Nevim co už dál / Old-czech style turkey breast with jasmine rice"
$ws.Range("C4").Value = "Jídlo 3 Jídlo Nevim co už dál / Grilled minced meat with roasted potatoes and spicy salad of roasted peppers"
$ws.Range("C5").Value = "Jídlo 4 Nevim co už dál/ Tagliolini with beef tenderloin sprinkled with Grana Padano Cheese"
$ws.Range("C6").Value = "Jídlo 5 Nevim co už dál  / Baked zander with vegetables in butter served with parsley potatoes"

$ws.Range("F2").Value = "Mlsná polévka / Potato soup"
$ws.Range("F3").Value = "Kuřecí vývarov / Lentil soup with sausages"
$ws.Range("F4").Value = "Porek / Beef consommé with meat and noodles"
$ws.Range("F5").Value = "Best polívka / Minestrone soup with pasta"
$ws.Range("F6").Value = "Pátková polívka / Bank holiday. We do not serve daily menu."

$ws.Range("A2").Value = 45859
$ws.Range("B2").Value = 45863

$ws.Range("F5").Select()
